$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.929.31'
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = '  +0.04%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.582.55'
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = '  +1.56%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.87'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  +1.48%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.31'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  +1.03%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("E8").Value = '  +2.63%  '
$ws.Range("E9").Value = '  +2.64%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.65'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  +3.04%  '
$ws.Range("E12").Value = '  -0.12%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '27.35'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = '  +0.43%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.045.28'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = '  +1.60%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '62.798.12'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = '  -0.05%  '
$ws.Range("E16").Value = '  +3.07%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.584.87'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = '  +2.11%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.30'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = '  -0.13%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '342.55'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = '  +2.01%  '
$ws.Range("E20").Value = '  +1.97%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.68'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  -0.84%  '
$ws.Range("E22").Value = '  -0.16%  '
$ws.Range("E23").Value = '  -1.04%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '67.35'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = '  +3.29%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.717.24'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = '  +1.68%  '
$ws.Range("E26").Value = '  -0.54%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.59'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = '  +0.45%  '
$ws.Range("E29").Value = '  +9.03%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.32'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = '  +0.04%  '
$ws.Range("E31").Value = '  -0.99%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.93'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = '  +4.46%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0₃0820'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = '  +1.43%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '465.19'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = '  +15.01%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '175.01'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.60'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = '  +4.54%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.403'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = '  +1.44%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.09'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = '  +0.04%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.54'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = '  +5.09%  '
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("E42").Value = '  -1.46%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '158.12'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = '  +4.77%  '
$ws.Range("E44").Value = '  +1.30%  '
$ws.Range("E45").Value = '  +6.68%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '21.23'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = '  +2.42%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0542'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = '  +1.76%  '
$ws.Range("E48").Value = '  +0.79%  '
$ws.Range("E49").Value = '  -0.27%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '18.43'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = '  +1.26%  '
$ws.Range("E51").Value = '  +1.06%  '
